# Weekly driver report update for 2025-04-20
# Refresh the "Bad Drivers" and "Good Drivers" tables on the Driver Summary
# sheet with this week's roaming-impact numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bad Drivers (rows 3-7) + Totals (row 8) ---------------------------
$ws.Range("C3").Value = 58
$ws.Range("D3").Value = 93.40000000000001

$ws.Range("A4").Value = 'Intel(R) Wi-Fi 6 AX201 160MHz - 23.30.0.6'
$ws.Range("B4").Value = 10
$ws.Range("C4").Value = 608
$ws.Range("D4").Value = 95.40000000000001

$ws.Range("A5").Value = 'Intel(R) Wi-Fi 6E AX210 160MHz - 23.80.0.7'
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 17
$ws.Range("D5").Value = 95.5

$ws.Range("A6").Value = 'Intel(R) Wi-Fi 6 AX201 160MHz - 23.80.1.3'
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 9
$ws.Range("D6").Value = 98.40000000000001

$ws.Range("A7").Value = 'Intel(R) Wi-Fi 6 AX201 160MHz - 22.20.0.6'
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 98.5

$ws.Range("B8").Value = 18
$ws.Range("C8").Value = 693

# --- Good Drivers (rows 16-29) ------------------------------------------
$ws.Range("A16").Value = 'Intel(R) Wi-Fi 6E AX210 160MHz - 22.0.1.5'
$ws.Range("B16").Value = 156943
$ws.Range("D16").Value = 100

$ws.Range("A17").Value = 'Intel(R) Wi-Fi 6E AX210 160MHz - 23.120.0.3'
$ws.Range("B17").Value = 34181
$ws.Range("D17").Value = 99.90000000000001
$ws.Range("E17").Value = "'2025-02-05"

$ws.Range("A18").Value = 'Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4'
$ws.Range("B18").Value = 445055
$ws.Range("D18").Value = 99.90000000000001
$ws.Range("E18").Value = "'2024-11-10"

$ws.Range("A19").Value = 'Intel(R) Wi-Fi 6E AX210 160MHz - 23.20.1.1'
$ws.Range("B19").Value = 13533
$ws.Range("D19").Value = 100
$ws.Range("E19").Value = "'2023-12-19"

$ws.Range("A20").Value = 'Intel(R) Wi-Fi 6E AX210 160MHz - 22.170.2.1'
$ws.Range("B20").Value = 19083
$ws.Range("E20").Value = "'2022-08-30"

$ws.Range("A21").Value = 'Intel(R) Wi-Fi 6E AX210 160MHz - 22.100.0.3'
$ws.Range("B21").Value = 12988
$ws.Range("D21").Value = 100
$ws.Range("E21").Value = "'2022-05-01"

$ws.Range("A22").Value = 'Intel(R) Wi-Fi 6E AX210 160MHz - 22.130.0.5'
$ws.Range("B22").Value = 18738
$ws.Range("D22").Value = 99.90000000000001
$ws.Range("E22").Value = "'2022-03-14"

$ws.Range("A23").Value = 'Intel(R) Wi-Fi 6E AX210 160MHz - 22.110.1.1'
$ws.Range("B23").Value = 42024
$ws.Range("E23").Value = "'2022-01-01"

$ws.Range("A24").Value = 'Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9'
$ws.Range("B24").Value = 77849
$ws.Range("D24").Value = 99.90000000000001
$ws.Range("E24").Value = "'2021-08-18"

$ws.Range("A25").Value = 'Intel(R) Wi-Fi 6E AX210 160MHz - 22.70.0.6'
$ws.Range("B25").Value = 15504
$ws.Range("E25").Value = "'2021-06-28"

$ws.Range("A26").Value = 'Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1'
$ws.Range("B26").Value = 34244
$ws.Range("D26").Value = 100
$ws.Range("E26").Value = "'2021-04-27"

$ws.Range("A27").Value = 'Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2'
$ws.Range("B27").Value = 59673
$ws.Range("E27").Value = "'2020-08-05"

$ws.Range("A28").Value = 'Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6'
$ws.Range("B28").Value = 113652
$ws.Range("E28").Value = "'2020-01-06"

$ws.Range("A29").Value = 'Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1'
$ws.Range("B29").Value = 56018
